# Insert a new "Sheet1" datasheet between "clients" and "vendors" that
# summarizes Total Expenses per client (Name), reproducing the author's
# "updated datatable output to sheet1 missing email" edit.

$wb = $excel.ActiveWorkbook

# Insert the new sheet right before "vendors" so the tab order becomes:
# clients, Sheet1, vendors, vendor_inventory, incoming_shipment, Test,
# expense_reports, out_of_stock
$ws = $wb.Worksheets.Add($wb.Worksheets.Item("vendors"))
$ws.Name = "Sheet1"

# Header row
$ws.Range("A1").Value = "Name"
$ws.Range("B1").Value = "Total Expenses"

# Data rows (Name + Total Expenses), matching the values from the
# author's commit (note: no Email column -- that's the "missing email").
$ws.Range("A2").Value = "Carson Goble`n"
$ws.Range("B2").Value = 59.465

$ws.Range("A3").Value = "Cayden Doyle`n`n"
$ws.Range("B3").Value = 1876.761

$ws.Range("A4").Value = "Aiden Herrera`n`n"
$ws.Range("B4").Value = 3009.186

$ws.Range("A5").Value = "Cindy"
$ws.Range("B5").Value = 3.98

# Wrap-text formatting down columns A and D through row 18 (matches the
# extra cellXfs style added for this sheet).
$ws.Range("A1:A18").WrapText = $true
$ws.Range("D1:D18").WrapText = $true

# Row heights for the wrapped multi-line name cells.
$ws.Rows.Item(2).RowHeight = 37.5
$ws.Rows.Item(3).RowHeight = 50
$ws.Rows.Item(4).RowHeight = 50

# Column widths (B wider for names, C sized for typical content).
$ws.Columns.Item(2).ColumnWidth = 22.25651041666667
$ws.Columns.Item(3).ColumnWidth = 13.619791666666666

# Final UI state: new sheet is the active/selected tab with C8 selected.
$null = $ws.Range("C8").Select()
